$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28, shifting existing rows 28-71 down to 29-72.
$ws.Rows("28:28").Insert()

# Populate the newly inserted row 28 with the new record.
$ws.Cells.Item(28, 1).Value = 11
$ws.Cells.Item(28, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(28, 3).Value = "Bíobío"
$ws.Cells.Item(28, 4).Value = 44467
$ws.Cells.Item(28, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(28, 5).Value = 8
$ws.Cells.Item(28, 6).Value = 100112043
$ws.Cells.Item(28, 7).Value = "Pepino ensalada"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 100
$ws.Cells.Item(28, 11).Value = 15000
$ws.Cells.Item(28, 12).Value = 16000
$ws.Cells.Item(28, 13).Value = 15500
$ws.Cells.Item(28, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(28, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(28, 16).Value = 258
$ws.Cells.Item(28, 17).Value = 60
$ws.Cells.Item(28, 18).Value = "Hortaliza"
